$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.587.45"
$ws.Range("E2").Value = "  +5.48%  "
$ws.Range("D3").Value = "3.647.55"
$ws.Range("E3").Value = "  +5.82%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "193.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.92%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.646"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.39%  "
$ws.Range("D8").Value = "3.641.74"
$ws.Range("E8").Value = "  +5.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.183"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.675"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.42%  "
$ws.Range("E13").Value = "  +5.84%  "
$ws.Range("E14").Value = "  +5.84%  "
$ws.Range("D15").Value = "4.228.24"
$ws.Range("E15").Value = "  +5.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.47%  "
$ws.Range("D17").Value = "3.645.62"
$ws.Range("E17").Value = "  +5.70%  "
$ws.Range("D18").Value = "70.659.44"
$ws.Range("E18").Value = "  +5.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.67%  "
$ws.Range("E20").Value = "  +3.15%  "
$ws.Range("E21").Value = "  +4.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "488.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "91.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.91%  "
$ws.Range("E27").Value = "  +6.94%  "
$ws.Range("E28").Value = "  +5.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.22%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.121"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "625.38"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.64%  "
$ws.Range("E34").Value = "  +4.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "65.86"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.92%  "
$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.411"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.73%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0823"
$ws.Range("E38").Value = "  +9.59%  "
$ws.Range("E39").Value = "  -1.66%  "
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("D42").Value = "3.299.11"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.35%  "
$ws.Range("E44").Value = "  +11.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0455"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.03"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.19%  "
$ws.Range("E47").Value = "  +2.10%  "
$ws.Range("E48").Value = "  +3.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +7.10%  "
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.29%  "
